# Add a new "Vị trí kiêm nhiệm" (concurrent position) column between the
# existing "Vị trí công việc" column (D) and "Mục tiêu cá nhân" column (E),
# pushing the latter (and its data) one column to the right (E -> F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; this shifts column E (and anything after it)
# one column to the right, carrying formatting/styles along with it.
$ws.Range("E1").EntireColumn.Insert()

# Give the newly inserted column E its header text (the data rows below
# stay empty, same as the source workbook).
$ws.Range("E1").Value = "Vị trí kiêm nhiệm"

# Match the author's final cursor position left in the worksheet.
$ws.Range("E8").Select() | Out-Null
